$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 54 with the new contract record (previously only B54 held a
# leftover "600362" label with no other data).
$ws.Range("A54").Value = "江西铜业"
$ws.Range("B54").Value = 600362
$ws.Range("C54").Value = 20201211
$ws.Range("D54").Value = 20210611
$ws.Range("E54").Value = 1000
$ws.Range("F54").Value = "未了结"
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 8.35
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 13.84
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = "000000028365"
$ws.Range("O54").Value = "融券卖出"
$ws.Range("P54").Value = 21040
$ws.Range("Q54").Value = 1000
$ws.Range("R54").Value = 21.04

# B54 previously carried formatting (style index 1) inherited from the old
# placeholder value; the new row uses the default (unstyled) cells, same as
# other plain numeric columns.
$ws.Range("B54").ClearFormats()

# Update the view: scroll back to the top of the sheet and move the active
# selection to E33 (was A19/F47).
$ws.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 1
